# Apply hybrid bold + color highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet
# paragraphs, per the commit's "quantitative metrics highlighting"
# feature. Each target paragraph is located by its exact original
# text, then the numeric substrings within it are located (in left
# to right order) and given Bold + RGB(44,62,80) == hex 2C3E50
# formatting, which Word stores as <w:b/><w:color w:val="2C3E50"/>.

$d = $word.ActiveDocument

# Word COM Font.Color wants a BGR-packed long; 0x2C3E50 (RGB) ->
# packed as 0x00503E2C.
$metricColor = 5258796

function Highlight-Metric($paragraphRange, [string]$needle) {
    $found = $paragraphRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $paragraphRange.Font.Bold = 1
        $paragraphRange.Font.Color = $metricColor
    }
    return $found
}

# Map of exact original paragraph text -> ordered list of numeric
# substrings (first occurrence of each, in left-to-right order) that
# should become bold + colored.
$targets = @(
    @{
        Text = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
        Metrics = @("23%", "64%")
    },
    @{
        Text = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
        Metrics = @("87%", "71%", "±4.2%", "±2.1%")
    },
    @{
        Text = "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
        Metrics = @("1,200")
    },
    @{
        Text = "• Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+"
        Metrics = @("`$400M", "`$1B")
    },
    @{
        Text = "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M"
        Metrics = @("73.5%", "`$4.7M")
    },
    @{
        Text = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
        Metrics = @("87%", "71%")
    }
)

$count = $d.Paragraphs.Count

foreach ($target in $targets) {
    $matched = $false
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $pText = $p.Range.Text
        # Paragraph Range.Text includes a trailing paragraph mark
        # (usually \r); trim any trailing CR/LF before comparing.
        $pText = $pText.TrimEnd("`r", "`n")
        if ($pText -eq $target.Text) {
            $matched = $true
            foreach ($metric in $target.Metrics) {
                $pr = $p.Range
                Highlight-Metric $pr $metric | Out-Null
            }
            break
        }
    }
    if (-not $matched) {
        Write-Host "WARNING: target paragraph not found: $($target.Text)"
    }
}

Write-Host "done"
